# Update the "datetimeFigureOut" Date Placeholder field text from 2/7/2021
# to 2/10/2021 across the slide master and every slide layout (the footer
# date stamp shown on every slide is inherited from these).
$p = $ppt.ActivePresentation

$oldDate = "2/7/2021"
$newDate = "2/10/2021"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout off the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $cl = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $cl.Shapes
}
